# Minor changes to comments
# Rebuild the APITestCaseData sheet: add a header row, re-letter the
# first column to use lower-camel-case test method names, and apply
# header/data styling (bold white-on-blue header, thin-bordered cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the shared-strings table gets rebuilt in
# the exact order the new layout needs.
$ws.Cells.Clear()

# --- Data rows (write these first so the "Weather_API_Test_NNN" ids and
# the description/message strings land at the same shared-string indexes
# as before) -----------------------------------------------------------
$ws.Range("B2").Value = "Weather_API_Test_001"
$ws.Range("B3").Value = "Weather_API_Test_002"
$ws.Range("B4").Value = "Weather_API_Test_003"
$ws.Range("B5").Value = "Weather_API_Test_004"
$ws.Range("B6").Value = "Weather_API_Test_005"

$ws.Range("C2").Value = "Get weather data in JSON format."
$ws.Range("D2").Value = "API response contains weather data for requested city in JSON format."
$ws.Range("E2").Value = "Unsuccessful API response."

$ws.Range("C3").Value = "Get weather data in XML format."
$ws.Range("D3").Value = "API response contains weather data for requested city in XML format."
$ws.Range("E3").Value = "Unsuccessful API response."

$ws.Range("C4").Value = "Get temperature data in imperial units."
$ws.Range("D4").Value = "API response contains temperature data in imperial units for requested city in JSON format."
$ws.Range("E4").Value = "Unsuccessful API response."

$ws.Range("C5").Value = "Get temperature data in metric units."
$ws.Range("D5").Value = "API response contains temperature data in metric units for requested city in JSON format."
$ws.Range("E5").Value = "Unsuccessful API response."

$ws.Range("C6").Value = "Get weather data in non-English language."
$ws.Range("D6").Value = "API response contains weather data for requested city in specified language."
$ws.Range("E6").Value = "Unsuccessful API response."

# --- Header row ---------------------------------------------------------
$ws.Range("A1").Value = "Test_Case_Name"
$ws.Range("B1").Value = "Test_Case_ID"
$ws.Range("C1").Value = "Test_Description"
$ws.Range("D1").Value = "Success_Message"
$ws.Range("E1").Value = "Failure_Message"

# --- Column A test-method names (lower camel case) ----------------------
$ws.Range("A2").Value = "getWeatherInJSONFormatTest"
$ws.Range("A3").Value = "getWeatherInXMLFormatTest"
$ws.Range("A4").Value = "getTempInFahrenheitUnitTest"
$ws.Range("A5").Value = "getTempInCelsiusUnitTest"
$ws.Range("A6").Value = "getWeatherInDiffLangTest"

# --- Styling -------------------------------------------------------------
# Build the two new styles on scratch cells off to the side, then copy
# just the formats onto the real ranges - this keeps the number of
# newly-materialised font/fill/border/cellXf entries to a minimum
# (direct per-cell property assignment would re-resolve a style on every
# single call).
$headerScratch = $ws.Range("H1")
$headerScratch.Font.ThemeColor = 2        # xlThemeColorLight1 -> theme="0" (white)
$headerScratch.Font.Bold = $true
$headerScratch.Interior.ThemeColor = 5    # xlThemeColorAccent1 -> theme="4" (blue)
$headerScratch.Borders.LineStyle = 1      # xlContinuous / thin
$headerScratch.HorizontalAlignment = -4108 # xlCenter
$headerScratch.VerticalAlignment = -4108   # xlCenter

$dataScratch = $ws.Range("H2")
$dataScratch.Borders.LineStyle = 1

$headerRange = $ws.Range("A1:E1")
$headerScratch.Copy()
$headerRange.PasteSpecial(-4122) # xlPasteFormats

$dataRange = $ws.Range("A2:E6")
$dataScratch.Copy()
$dataRange.PasteSpecial(-4122) # xlPasteFormats

$headerScratch.Clear()
$dataScratch.Clear()

# --- Column widths (best-fit for the new/longer content) -----------------
$ws.Columns.Item(3).ColumnWidth = 39.28515625
$ws.Columns.Item(4).ColumnWidth = 84.140625
$ws.Columns.Item(5).ColumnWidth = 25.5703125

# --- Selection mirrors the author's saved cursor position ----------------
[void]$ws.Range("A7").Select()

Write-Host "done"
